# Insert a new data row into the "Coco" sheet at row 18 (pushing existing
# rows 18-67 down to 19-68) and populate it with the new weekly price entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 18. This shifts rows 18..67 down
# to 19..68 and duplicates the formatting of the row above (row 17),
# matching the existing data rows' look (style carried on column D).
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new record's values.
$ws.Cells.Item(18, 1).Value = 10
$ws.Cells.Item(18, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(18, 3).Value = "La Araucanía"
$ws.Cells.Item(18, 4).NumberFormat = $ws.Cells.Item(19, 4).NumberFormat
$ws.Cells.Item(18, 4).Value = 44715
$ws.Cells.Item(18, 5).Value = 9
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100108
$ws.Cells.Item(18, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(18, 9).Value = 100108007
$ws.Cells.Item(18, 10).Value = "Coco"
$ws.Cells.Item(18, 11).Value = "Sin especificar"
$ws.Cells.Item(18, 12).Value = "Primera"
$ws.Cells.Item(18, 13).Value = 10
$ws.Cells.Item(18, 14).Value = 30000
$ws.Cells.Item(18, 15).Value = 30000
$ws.Cells.Item(18, 16).Value = 30000
$ws.Cells.Item(18, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(18, 18).Value = "Perú"
$ws.Cells.Item(18, 19).Value = 1500
$ws.Cells.Item(18, 20).Value = 20
